# Apply the "noise" update to the generated data sheet.
# - Columns K (11) and L (12) get slightly narrower widths.
# - Columns I, J, K, L (and a few N cells) get new noisy values
#   for rows 1-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (K and L got 1 px / ~0.33 char narrower) ---
$ws.Columns.Item(11).ColumnWidth = 13.72
$ws.Columns.Item(12).ColumnWidth = 14.72

# --- Row 1 ---
$ws.Range("I1").Value = 223
$ws.Range("J1").Value = 89.243538400000006
$ws.Range("K1").Value = 0.0011452752026068769
$ws.Range("L1").Value = 0.00017670182326759047

# --- Row 2 ---
$ws.Range("I2").Value = 273
$ws.Range("J2").Value = 198.4210301
$ws.Range("K2").Value = 0.0013440120782932663
$ws.Range("L2").Value = 0.00018336590789271059

# --- Row 3 ---
$ws.Range("I3").Value = 272
$ws.Range("J3").Value = 241.57603230000001
$ws.Range("K3").Value = 0.0012797381774800609
$ws.Range("L3").Value = 0.0001934054255676059

# --- Row 4 ---
$ws.Range("I4").Value = 277
$ws.Range("J4").Value = 285.82469040000001
$ws.Range("K4").Value = 0.001428688523450905
$ws.Range("L4").Value = 0.00020959038114722305

# --- Row 5 ---
$ws.Range("I5").Value = 241
$ws.Range("J5").Value = 189.54301290000001
$ws.Range("K5").Value = 0.0016855521033674048
$ws.Range("L5").Value = 0.00019848763714879484

# --- Row 6 ---
$ws.Range("I6").Value = 788
$ws.Range("J6").Value = 809.55318469999997
$ws.Range("K6").Value = 0.00156134978185829
$ws.Range("L6").Value = 0.00020824307857343987

# --- Row 7 ---
$ws.Range("I7").Value = 1286
$ws.Range("J7").Value = 1712.7446198
$ws.Range("K7").Value = 0.0018178502689683018
$ws.Range("L7").Value = 0.0001805025436016349

# --- Row 8 ---
$ws.Range("I8").Value = 288
$ws.Range("J8").Value = 215.7881635
$ws.Range("K8").Value = 0.0012400114727397771
$ws.Range("L8").Value = 0.00022390849220586047
$ws.Range("N8").Value = 62

# --- Row 9 ---
$ws.Range("I9").Value = 314
$ws.Range("J9").Value = 234.1855151
$ws.Range("K9").Value = 0.0013339202775304315
$ws.Range("L9").Value = 0.00017829159837615108
$ws.Range("N9").Value = 83

# --- Row 10 ---
$ws.Range("I10").Value = 375
$ws.Range("J10").Value = 817.5757691
$ws.Range("K10").Value = 0.0011519689145058631
$ws.Range("L10").Value = 0.00020460594234797855
$ws.Range("N10").Value = 43

# --- Row 11 ---
$ws.Range("I11").Value = 273
$ws.Range("J11").Value = 323.64125680000001
$ws.Range("K11").Value = 0.0014627231347541514
$ws.Range("L11").Value = 0.00017274467533679911

# --- Row 12 ---
$ws.Range("I12").Value = 273
$ws.Range("J12").Value = 313.73481399999997
$ws.Range("K12").Value = 0.0013883434331662947
$ws.Range("L12").Value = 0.0001514347763090829

# --- Row 13 ---
$ws.Range("I13").Value = 272
$ws.Range("J13").Value = 472.09343189999998
$ws.Range("K13").Value = 0.0014953237659585117
$ws.Range("L13").Value = 0.00017036604357627353

Write-Host "edit applied"
